$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Table row properties: every row gets <w:cantSplit/>, the header
#        row additionally gets <w:tblHeader/>. -----------------------------
$t.Rows.AllowBreakAcrossPages = $false
$t.Rows.First.HeadingFormat = $true

# --- 2. "General" row about a native English speaker review: fill in the
#        (until now empty) Committee response cell. ------------------------
$t.Cell(6, 3).Range.Text = "Native English speaker review took place. "

# --- 3. Annex C row ("Suggestion: Background and context ..."): fill in
#        the (until now empty) Committee response cell. --------------------
$t.Cell(23, 3).Range.Text = "Suggestion adopted. "

# --- 4. Pagination shifted by one row because of the text added above, so
#        the rendered page break marker moves from the "FAT 286" cell to
#        the "Annex C" cell. ------------------------------------------------
$annexCCell = $t.Cell(23, 1).Range
$annexCXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7FB1BEC4" w14:textId="349879FB" w:rsidR="0000436A" w:rsidRDefault="00891288" w:rsidP="00856943"><w:pPr><w:tabs><w:tab w:val="left" w:pos="4780"/></w:tabs></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Annex C</w:t></w:r></w:p>'
$annexCCell.InsertXML($annexCXml) | Out-Null

$fat286Cell = $t.Cell(24, 1).Range
$fat286Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="01D66D54" w14:textId="701A06EA" w:rsidR="003E151F" w:rsidRDefault="003E151F" w:rsidP="00856943"><w:pPr><w:tabs><w:tab w:val="left" w:pos="4780"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">FAT 286 </w:t></w:r></w:p>'
$fat286Cell.InsertXML($fat286Xml) | Out-Null

# --- 5. FAT 286 row: fill in the (until now empty) Committee response
#        cell explaining the bibliography fix. ------------------------------
$t.Cell(24, 3).Range.Text = "“286” is the number of the VDA standard χMCF 3.0, which was referenced as [13]. " + [char]11 + "Bibliography item extended, accordingly. "
